$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 47, shifting the existing data (old rows 47-144)
# down to rows 49-146.
$ws.Range("A47:A48").EntireRow.Insert()

# New row 47: Haba, "Primera" quality, Region Metropolitana, dated 2021-10-13 (44482)
$ws.Cells.Item(47, 1).Value = 6
$ws.Cells.Item(47, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(47, 3).Value = "Metropolitana"
$ws.Cells.Item(47, 4).Value = 44482
$ws.Cells.Item(47, 5).Value = 13
$ws.Cells.Item(47, 6).Value = 100112026
$ws.Cells.Item(47, 7).Value = "Haba"
$ws.Cells.Item(47, 8).Value = "Sin especificar"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 900
$ws.Cells.Item(47, 11).Value = 5000
$ws.Cells.Item(47, 12).Value = 6000
$ws.Cells.Item(47, 13).Value = 5556
$ws.Cells.Item(47, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(47, 15).Value = "Región Metropolitana"
$ws.Cells.Item(47, 16).Value = 222
$ws.Cells.Item(47, 17).Value = 25
$ws.Cells.Item(47, 18).Value = "Hortaliza"

# New row 48: Haba, "Segunda" quality, Region Metropolitana, dated 2021-10-13 (44482)
$ws.Cells.Item(48, 1).Value = 6
$ws.Cells.Item(48, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(48, 3).Value = "Metropolitana"
$ws.Cells.Item(48, 4).Value = 44482
$ws.Cells.Item(48, 5).Value = 13
$ws.Cells.Item(48, 6).Value = 100112026
$ws.Cells.Item(48, 7).Value = "Haba"
$ws.Cells.Item(48, 8).Value = "Sin especificar"
$ws.Cells.Item(48, 9).Value = "Segunda"
$ws.Cells.Item(48, 10).Value = 300
$ws.Cells.Item(48, 11).Value = 4000
$ws.Cells.Item(48, 12).Value = 4000
$ws.Cells.Item(48, 13).Value = 4000
$ws.Cells.Item(48, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(48, 15).Value = "Región Metropolitana"
$ws.Cells.Item(48, 16).Value = 160
$ws.Cells.Item(48, 17).Value = 25
$ws.Cells.Item(48, 18).Value = "Hortaliza"
